$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 203
$ws.Range("F3").Value = 5508
$ws.Range("F6").Value = 30
$ws.Range("F7").Value = 650
$ws.Range("F8").Value = 632
$ws.Range("F9").Value = 1072
$ws.Range("F11").Value = 1530
$ws.Range("F12").Value = 5034
$ws.Range("F14").Value = 229
$ws.Range("F15").Value = 199
$ws.Range("F16").Value = 10
$ws.Range("F18").Value = 4298
$ws.Range("F19").Value = 199
$ws.Range("F20").Value = 1144
$ws.Range("F21").Value = 117
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 159
$ws.Range("F26").Value = 57
$ws.Range("F29").Value = 339
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 203
$ws.Range("F4").Value = 5508
$ws.Range("F7").Value = 30
$ws.Range("F8").Value = 650
$ws.Range("F9").Value = 632
$ws.Range("F10").Value = 1072
$ws.Range("F12").Value = 1530
$ws.Range("F13").Value = 5034
$ws.Range("F15").Value = 229
$ws.Range("F16").Value = 199
$ws.Range("F17").Value = 10
$ws.Range("F19").Value = 4298
$ws.Range("F20").Value = 199
$ws.Range("F21").Value = 1144
$ws.Range("F22").Value = 117
$ws.Range("F25").Value = 53
$ws.Range("F26").Value = 159
$ws.Range("F27").Value = 57
$ws.Range("F30").Value = 339

$wb.Save()
